$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 holds the item data:
#  H7 (merged H7:K7) = current balance ("الرصيد الحالي"): "9:0" -> "8:0"
#  P7 = selling price ("سعر البيع"): "40.0000" -> "80.0000"
#  Q7 = number of transactions ("عدد التعاملات"): "1:0" -> "2:0"

$ws.Range("H7").Value = "8:0"

# P7 is formatted with a numeric number format (0.00), so a plain text-looking
# assignment would be auto-converted to the number 80. Temporarily switch the
# cell to a text format, write the literal string, then restore the original
# number format so the value is stored as text (matching the source data).
$cell = $ws.Range("P7")
$oldFormat = $cell.NumberFormat
$cell.NumberFormat = "@"
$cell.Value = "80.0000"
$cell.NumberFormat = $oldFormat

$ws.Range("Q7").Value = "2:0"
